$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set new values for country_name (J2) and region_name (K2)
$ws.Range("J2").Value = "Germany"
$ws.Range("K2").Value = "Bavaria"

# Update the selected cell/range on the sheet view
$ws.Range("K7").Select()

# Update the workbook window position/size (best-effort; host may only
# persist a subset of these window-chrome properties to the OOXML)
$win = $excel.ActiveWindow
$win.Left = 33300
$win.Top = -1060
$win.Width = 30240
$win.Height = 17800
